# Applies the EOL-RIR target price Max sensitivity update:
#  - Column C header year changes from 2020 to 2030 on every sheet
#  - Recomputed probability values for columns B:E on every sheet
$wb = $excel.ActiveWorkbook

# Sheet 1: Neodymium
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(1, 3).Value = 2030
$ws.Cells.Item(2, 3).Value = [double]"2.438729849902814E-05"
$ws.Cells.Item(2, 4).Value = [double]"0.1236014393403301"
$ws.Cells.Item(2, 5).Value = [double]"0.584289849462961"
$ws.Cells.Item(3, 2).Value = [double]"2.183980475909259E-12"
$ws.Cells.Item(3, 3).Value = [double]"0.001187950236302112"
$ws.Cells.Item(3, 4).Value = [double]"0.1155224455664189"
$ws.Cells.Item(3, 5).Value = [double]"0.4948953125759609"
$ws.Cells.Item(4, 2).Value = [double]"3.409259119931335E-14"
$ws.Cells.Item(4, 3).Value = [double]"0.001073104153047114"
$ws.Cells.Item(4, 4).Value = [double]"0.08248514375558298"
$ws.Cells.Item(4, 5).Value = [double]"0.4143225657740967"
$ws.Cells.Item(5, 3).Value = [double]"2.380640363208162E-08"
$ws.Cells.Item(5, 4).Value = [double]"0.004220506989166603"
$ws.Cells.Item(5, 5).Value = [double]"0.03292316497354641"

# Sheet 2: Dysprosium
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(1, 3).Value = 2030
$ws.Cells.Item(2, 3).Value = [double]"2.763358668632657E-05"
$ws.Cells.Item(2, 4).Value = [double]"0.1067081964413131"
$ws.Cells.Item(2, 5).Value = [double]"0.6620669446318078"
$ws.Cells.Item(3, 3).Value = [double]"0.001346082914235243"
$ws.Cells.Item(3, 4).Value = [double]"0.0997334001988444"
$ws.Cells.Item(3, 5).Value = [double]"0.5607727531308762"
$ws.Cells.Item(4, 3).Value = [double]"0.001215949222004488"
$ws.Cells.Item(4, 4).Value = [double]"0.07121147593698561"
$ws.Cells.Item(4, 5).Value = [double]"0.4694746545164075"
$ws.Cells.Item(5, 3).Value = [double]"2.697536664354322E-08"
$ws.Cells.Item(5, 4).Value = [double]"0.003643668644034803"
$ws.Cells.Item(5, 5).Value = [double]"0.03730569555791484"

# Sheet 3: Copper
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(1, 3).Value = 2030
$ws.Cells.Item(2, 2).Value = [double]"3.278497091721097E-06"
$ws.Cells.Item(2, 3).Value = [double]"0.003050246220774824"
$ws.Cells.Item(2, 4).Value = [double]"0.8588525024148915"
$ws.Cells.Item(2, 5).Value = [double]"0.9697180417558001"
$ws.Cells.Item(3, 2).Value = [double]"2.229370101113288E-05"
$ws.Cells.Item(3, 3).Value = [double]"0.01103680953317707"
$ws.Cells.Item(3, 4).Value = [double]"0.6168778258139214"
$ws.Cells.Item(3, 5).Value = [double]"0.6805195557530022"
$ws.Cells.Item(4, 2).Value = [double]"6.612099022439717E-05"
$ws.Cells.Item(4, 3).Value = [double]"0.002936712507453067"
$ws.Cells.Item(4, 4).Value = [double]"0.4416445656454649"
$ws.Cells.Item(4, 5).Value = [double]"0.6055420661490584"
$ws.Cells.Item(5, 2).Value = [double]"2.076994439830034E-05"
$ws.Cells.Item(5, 3).Value = [double]"0.006488512101692222"
$ws.Cells.Item(5, 4).Value = [double]"0.8181956015862759"
$ws.Cells.Item(5, 5).Value = [double]"0.7122850722825392"

# Sheet 4: Raw silicon
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(1, 3).Value = 2030
$ws.Cells.Item(2, 2).Value = [double]"4.96652837099915E-07"
$ws.Cells.Item(2, 3).Value = [double]"0.0005194103187024707"
$ws.Cells.Item(2, 4).Value = [double]"0.4783497257479737"
$ws.Cells.Item(2, 5).Value = [double]"1.251484198240786"
$ws.Cells.Item(3, 2).Value = [double]"5.30035999530297E-07"
$ws.Cells.Item(3, 3).Value = [double]"0.001746085699123397"
$ws.Cells.Item(3, 4).Value = [double]"0.2315523203444246"
$ws.Cells.Item(3, 5).Value = [double]"0.5836018294026056"
$ws.Cells.Item(4, 2).Value = [double]"3.397047964529607E-06"
$ws.Cells.Item(4, 3).Value = [double]"0.0004869161330086471"
$ws.Cells.Item(4, 4).Value = [double]"0.2164427902680356"
$ws.Cells.Item(4, 5).Value = [double]"0.6322978441499242"
$ws.Cells.Item(5, 2).Value = [double]"1.823860200208514E-06"
$ws.Cells.Item(5, 3).Value = [double]"0.0006186310519116159"
$ws.Cells.Item(5, 4).Value = [double]"0.4280362603778582"
$ws.Cells.Item(5, 5).Value = [double]"0.8793987462358255"
